# Folder 01 - 02nd page updated
# Adds 48 new Brahmi word-translation entries (rows 38-85) to the
# "Folder 01 Translations" worksheet, matching the new dictionary data
# appended to the shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Set cell values for new rows 38-85 ---
$ws.Range('A38').Value = 'adiya'
$ws.Range('D38').Value = 'a+(da+i)+ya'

$ws.Range('A39').Value = 'ananika'
$ws.Range('B39').Value = 'an irrigation expert'
$ws.Range('D39').Value = 'a+na+(na+i)+ka'

$ws.Range('A40').Value = 'anikata'
$ws.Range('B40').Value = 'a body guard'
$ws.Range('D40').Value = 'a+(na+i)+ka+ta'

$ws.Range('A41').Value = 'ati'
$ws.Range('B41').Value = 'elephnat'
$ws.Range('D41').Value = 'a+(ta+i)'

$ws.Range('A42').Value = 'ati-acariya'
$ws.Range('B42').Value = 'leader of a company of a elephnat trainers'
$ws.Range('D42').Value = 'a+(ta+i)+a+ca+(ra+i)+ya'

$ws.Range('A43').Value = 'ati-ajariya'
$ws.Range('B43').Value = 'trainer of elephants'
$ws.Range('D43').Value = 'a+(ta+i)+a+ja+(ra+i)+ya'

$ws.Range('A44').Value = 'ati-adika'
$ws.Range('B44').Value = 'superintendent of elephnats'
$ws.Range('D44').Value = 'a+(ta+i)+a+(da+i)+ka'

$ws.Range('A45').Value = 'Ati-matakaha'
$ws.Range('B45').Value = 'one whose mother is of the Atri or Atreya-gotra'
$ws.Range('D45').Value = 'A+(ta+i)+ma+ta+ka+ha'

$ws.Range('A46').Value = 'atireka'
$ws.Range('B46').Value = 'remainder,residue'
$ws.Range('D46').Value = 'a+(ta+i)+(ra+e)+ka'

$ws.Range('A47').Value = 'ativas''ika'
$ws.Range('B47').Value = 'pupil,disciple'
$ws.Range('D47').Value = 'a+(ta+i)+va+(s''a+i)+ka'

$ws.Range('A48').Value = 'atevas''ika'
$ws.Range('D48').Value = 'a+(ta+e)+va+(s''a+i)+ka'

$ws.Range('A49').Value = 'atevahika'
$ws.Range('D49').Value = 'a+(ta+e)+va+(ha+i)+ka'

$ws.Range('A50').Value = 'adaka'
$ws.Range('B50').Value = 'superintendent'
$ws.Range('D50').Value = 'a+da+ka'

$ws.Range('A51').Value = 'Ada-kacaka'
$ws.Range('B51').Value = 'place name prefixed to a personal name'
$ws.Range('D51').Value = 'A+da+ka+ca+ka'

$ws.Range('A52').Value = 'adakaha'
$ws.Range('D52').Value = 'a+da+ka+ha'

$ws.Range('A53').Value = 'adika'
$ws.Range('D53').Value = 'a+(da+i)+ka'

$ws.Range('A54').Value = 'Adiliya'
$ws.Range('D54').Value = 'A+(da+i)+(la+i)+ya'

$ws.Range('A55').Value = 'adeka'
$ws.Range('D55').Value = 'a+(da+e)+ka'

$ws.Range('A56').Value = 'anagata'
$ws.Range('B56').Value = 'future'
$ws.Range('D56').Value = 'a+na+ga+ta'

$ws.Range('A57').Value = 'Anada'
$ws.Range('D57').Value = 'A+na+da'

$ws.Range('A58').Value = 'anu jete'
$ws.Range('B58').Value = 'deputy of the alderman'
$ws.Range('D58').Value = 'a+(na+u)+(ja+e)+(ta+e)'

$ws.Range('A59').Value = 'Anudi'
$ws.Range('B59').Value = 'growing,increasing'
$ws.Range('D59').Value = 'A+(na+u)+(da+i)'

$ws.Range('A60').Value = 'Anudi-gamasi'
$ws.Range('D60').Value = 'A+(na+u)+(da+i)+ga+ma+(sa+i)'

$ws.Range('A61').Value = 'anubuti'
$ws.Range('B61').Value = 'experience'
$ws.Range('D61').Value = 'a+(na+u)+(ba+u)+(ta+i)'

$ws.Range('A62').Value = 'Anurada'
$ws.Range('B62').Value = 'name of an asterism'
$ws.Range('D62').Value = 'A+(na+u)+ra+da'

$ws.Range('A63').Value = 'Anuradi'
$ws.Range('B63').Value = 'feminine of Anurada'
$ws.Range('D63').Value = 'A+(na+u)+ra+(da+i)'

$ws.Range('A64').Value = 'Anuridi'
$ws.Range('D64').Value = 'A+(na+u)+(ra+i)+(da+i)'

$ws.Range('A65').Value = 'Anula'
$ws.Range('B65').Value = 'name of a person'
$ws.Range('D65').Value = 'A+(na+u)+la'

$ws.Range('A66').Value = 'Anulapi'
$ws.Range('B66').Value = 'tank'
$ws.Range('D66').Value = 'A+(na+u)+la+(pa+i)'

$ws.Range('A67').Value = 'Anulaya'
$ws.Range('D67').Value = 'A+(na+u)+la+ya'

$ws.Range('A68').Value = 'Anotata'
$ws.Range('B68').Value = 'name of a lake in the Himalayas'
$ws.Range('D68').Value = 'A+(na+o)+ta+ta'

$ws.Range('A69').Value = 'Anodi'
$ws.Range('B69').Value = 'faultless'
$ws.Range('D69').Value = 'A+(na+o)+(da+i)'

$ws.Range('A70').Value = 'Apaya'
$ws.Range('B70').Value = 'fearless'
$ws.Range('D70').Value = 'A+pa+ya'

$ws.Range('A71').Value = 'apara'
$ws.Range('B71').Value = 'another'
$ws.Range('D71').Value = 'a+pa+ra'

$ws.Range('A72').Value = 'aparimita'
$ws.Range('B72').Value = 'boundless'
$ws.Range('D72').Value = 'a+pa+(ra+i)+(ma+i)+ta'

$ws.Range('A73').Value = 'api'
$ws.Range('B73').Value = 'as a second member of ..'
$ws.Range('D73').Value = 'a+(pa+i)'

$ws.Range('A74').Value = 'Aba'
$ws.Range('B74').Value = 'fearless'
$ws.Range('D74').Value = 'A+ba'

$ws.Range('A75').Value = 'aba'
$ws.Range('B75').Value = 'mango'
$ws.Range('D75').Value = 'a+ba'

$ws.Range('A76').Value = 'Aba-adi'
$ws.Range('B76').Value = 'mango tree channel'
$ws.Range('D76').Value = 'A+ba+a+(da+i)'

$ws.Range('A77').Value = 'abaka'
$ws.Range('B77').Value = 'wife'
$ws.Range('D77').Value = 'a+ba+ka'

$ws.Range('A78').Value = 'Aba-tota'
$ws.Range('B78').Value = 'ferry named after a mango tree'
$ws.Range('D78').Value = 'A+ba+(ta+o)+ta'

$ws.Range('A79').Value = 'Aba-nagara''na'
$ws.Range('B79').Value = 'residents or citizens or mebers of Aba-nagara'
$ws.Range('D79').Value = 'A+ba+na+ga+ra''+na'

$ws.Range('A80').Value = 'Aba-nagariyana'
$ws.Range('D80').Value = 'A+ba+na+ga+(ra+i)+ya+na'

$ws.Range('A81').Value = 'Abaya'
$ws.Range('B81').Value = 'Fearless'
$ws.Range('D81').Value = 'A+ba+ya'

$ws.Range('A82').Value = 'abala'
$ws.Range('B82').Value = 'sourness'
$ws.Range('D82').Value = 'a+ba+la'

$ws.Range('A83').Value = 'Aba-velaka'
$ws.Range('B83').Value = 'Mango-tree-field'
$ws.Range('D83').Value = 'A+ba+(va+e)+la+ka'

$ws.Range('A84').Value = 'abi'
$ws.Range('B84').Value = 'an honoric title attached to names of princesses'
$ws.Range('D84').Value = 'a+(ba+i)'

$ws.Range('A85').Value = 'Abijhatiya'
$ws.Range('B85').Value = 'high-born'
$ws.Range('D85').Value = 'A+(ba+i)+jha+(ta+i)+ya'


# --- Apply cell formatting (style "2": 10pt Arial / theme text color) ---
# Column A (rows 38-85) and column D (rows 65-85) use the same style as
# the rest of column A in this sheet; copy it across so the new cells
# match the existing look-and-feel.
$ws.Range('A2').Copy()
$ws.Range('A38:A85').PasteSpecial(-4122)
$ws.Range('D65:D85').PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row heights for the newly added rows ---
$ws.Range('38:85').RowHeight = 15.75

# --- Column width adjustments ---
# Target OOXML widths are 54.140625 (col B) and 13.7109375 (col C); the
# ColumnWidth property only accepts character-width units quantized to
# 1/6, so these are the closest achievable values (54.1666... / 13.6666...).
$ws.Columns.Item(2).ColumnWidth = 53.25
$ws.Columns.Item(3).ColumnWidth = 12.75

# --- Selection / scroll position ---
$ws.Range('D73').Select()
$excel.ActiveWindow.ScrollRow = 63
